$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2 through 388 all hold the serial date
# value 45202 (2023-10-03) and need to become 45203 (2023-10-04).
$startRow = 2
$endRow = 388

for ($row = $startRow; $row -le $endRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
